$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LANDING GEARS")

$ws.Range("C2").Value = 2233.0999999999995

$ws.Range("D6").Value = -24.81304016837578
$ws.Range("D7").Value = 5.2796560834714334
$ws.Range("D8").Value = 20.460346603376504
$ws.Range("D9").Value = 3.8018897496753645
